# Update "想去人数" (want-to-go count) figures to the freshly regenerated
# values, as produced by the gh-pages output generation at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 11818
$ws1.Range("F8").Value  = 4418
$ws1.Range("F17").Value = 5132
$ws1.Range("F21").Value = 11372
$ws1.Range("F22").Value = 11335

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 11818
$ws4.Range("F8").Value  = 4418
$ws4.Range("F18").Value = 5132
$ws4.Range("F22").Value = 11372
$ws4.Range("F23").Value = 11335
